# Product Register bug corrections
# Highlight (yellow) the "Aviso Legal" block of 5 paragraphs:
#   8 - Aviso Legal (deve ser fixo, ...)
#   -  Destinado para Maiores de 18 Anos
#   -  Imagens Meramente Ilustrativas
#   -  Medidas Aproximadas, Podendo Haver Pequenas Variações
#   -  Guarde produto separado de outros, ...
#
# The first four paragraphs get the yellow highlight applied to every run
# AND to the paragraph mark itself (pPr/rPr) - i.e. the whole paragraph,
# pilcrow included. The last paragraph only gets its runs highlighted
# (paragraph mark left untouched), matching how a manual click-drag
# highlight selection that ends mid-paragraph would behave in Word.

$wdYellow = 7

$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $rng.Find.Found) {
        throw "Text not found: $needle"
    }
    return $rng.Paragraphs(1)
}

# Paragraphs 1-4: highlight runs + paragraph mark (use Range.Font so the
# pPr/rPr picks up the highlight too).
$targets = @(
    "8 - Aviso Legal",
    "Destinado para Maiores de 18 Anos",
    "Imagens Meramente Ilustrativas",
    "Medidas Aproximadas, Podendo Haver Pequenas Varia"
)

foreach ($needle in $targets) {
    $p = Get-ParagraphByText $needle
    $p.Range.Font.HighlightColorIndex = $wdYellow
}

# Paragraph 5: highlight runs only, leave the paragraph mark as-is.
$pLast = Get-ParagraphByText "Guarde produto separado de outros"
$pLast.Range.HighlightColorIndex = $wdYellow
